# exception tag update module implemented and tested
#
# Source sheet ("Sheet1") is a list of exceptions ending at row 49 (row 49
# is already a blank spacer row, height 13.8). This change appends two more
# blank spacer rows (50, 51) with the same row height, narrows column B
# (the long "Exception" text column) from ~75 chars to ~58 chars, and moves
# the active selection down to the new last row (B51).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- extend the sheet with two more blank rows, matching row 49's height ---
$ws.Rows(50).RowHeight = 13.8
$ws.Rows(51).RowHeight = 13.8

# --- narrow column B (the "Exception" text column) ---
$ws.Columns(2).ColumnWidth = 57.33

# --- move the view/selection to the new bottom of the sheet ---
[void]$ws.Range("A47").Select()
[void]$ws.Range("B51").Select()
